# Update latest output (run 142)

$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule": update Cost ($) and Unit Cost ($/ML) for rows 4-5 ---
$scheduleWs = $wb.Worksheets.Item("Schedule")

$scheduleWs.Range("E4").Value = 718.1740800000001
$scheduleWs.Range("F4").Value = 31.66552380952382

$scheduleWs.Range("E5").Value = 235.1175645
$scheduleWs.Range("F5").Value = 6.911157098765432

# --- Sheet "Detailed": update Price (B) and Type (C) columns ---
$detailedWs = $wb.Worksheets.Item("Detailed")

$detailedWs.Range("B37").Value = 32.72201

$detailedWs.Range("B38").Value = 59.44775

$detailedWs.Range("B39").Value = 62.35405
$detailedWs.Range("C39").Value = "historical"

$detailedWs.Range("B40").Value = 69.06473
$detailedWs.Range("C40").Value = "historical"

$detailedWs.Range("B41").Value = 73.2
$detailedWs.Range("C41").Value = "historical"

$detailedWs.Range("B42").Value = 73.2
$detailedWs.Range("C42").Value = "historical"

$detailedWs.Range("B43").Value = 64.89
$detailedWs.Range("C43").Value = "historical"

$detailedWs.Range("B44").Value = 62.91923
$detailedWs.Range("C44").Value = "historical"

$detailedWs.Range("B45").Value = 62.59197
$detailedWs.Range("C45").Value = "historical"

$detailedWs.Range("B46").Value = 57.06
$detailedWs.Range("C46").Value = "historical"

$detailedWs.Range("B47").Value = 67.88919
$detailedWs.Range("C47").Value = "historical"

$detailedWs.Range("B48").Value = 73.19
$detailedWs.Range("C48").Value = "historical"

$detailedWs.Range("B49").Value = 65

$detailedWs.Range("B50").Value = 60.34232

$detailedWs.Range("B51").Value = 59.43612

$detailedWs.Range("B54").Value = 57.06

$detailedWs.Range("B55").Value = 56.97997

$detailedWs.Range("B56").Value = 57.06

$detailedWs.Range("B59").Value = 57.06007

$detailedWs.Range("B61").Value = 63.58873

$detailedWs.Range("B62").Value = 64.91909

$detailedWs.Range("B63").Value = 57.31

$detailedWs.Range("B64").Value = 36.0601

$detailedWs.Range("B66").Value = 27.79781

$detailedWs.Range("B67").Value = 10.50351

$detailedWs.Range("B68").Value = 20.43807

$detailedWs.Range("B69").Value = 9.795540000000001

$detailedWs.Range("B70").Value = 21.78921

$detailedWs.Range("B71").Value = 20.69799

$detailedWs.Range("B72").Value = 9.293760000000001

$detailedWs.Range("B73").Value = 9.010070000000001

$detailedWs.Range("B74").Value = 0.70268

$detailedWs.Range("B75").Value = 8.898149999999999

$detailedWs.Range("B76").Value = 8.67529

$detailedWs.Range("B77").Value = 8.65668

$detailedWs.Range("B78").Value = 0.51

$detailedWs.Range("B79").Value = 8.643750000000001

$detailedWs.Range("B80").Value = 8.832649999999999

$detailedWs.Range("B81").Value = 8.770960000000001

$detailedWs.Range("B82").Value = 8.95082

$detailedWs.Range("B83").Value = -1.0529

$detailedWs.Range("B84").Value = 0.00003

$detailedWs.Range("B86").Value = 57.06007

$detailedWs.Range("B87").Value = 73.19

$detailedWs.Range("B88").Value = 91.6481

$detailedWs.Range("B89").Value = 96.67813

$detailedWs.Range("B90").Value = 105.79

$detailedWs.Range("B91").Value = 92.00064999999999

$detailedWs.Range("B92").Value = 78

$detailedWs.Range("B93").Value = 73.19

$detailedWs.Range("B97").Value = 73.20007
